$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
$ws.Range("A8").Characters(21, 2).Text = "46"
$ws.Range("C9").Characters(47, 10).Text = "11/20/2022"
$ws.Range("C9").Characters(27, 9).Text = "11/14/2022"

# --- Cells changing from numeric to text placeholder (copy style+value from a cell already in that state) ---
$ws.Range("D15").Copy($ws.Range("C15"))
$ws.Range("D15").Copy($ws.Range("C26"))
$ws.Range("D15").Copy($ws.Range("D27"))
$ws.Range("D15").Copy($ws.Range("F30"))
$ws.Range("E15").Copy($ws.Range("E27"))

# --- Cells changing from text placeholder to numeric (copy numeric style first, then set value) ---
$ws.Range("C16").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1
$ws.Range("C16").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 2
$ws.Range("C16").Copy($ws.Range("G22"))
$ws.Range("G22").Value = 2
$ws.Range("C16").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("C16").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 2
$ws.Range("C16").Copy($ws.Range("F28"))
$ws.Range("F28").Value = 2
$ws.Range("C16").Copy($ws.Range("C29"))
$ws.Range("C29").Value = 2
$ws.Range("C16").Copy($ws.Range("F29"))
$ws.Range("F29").Value = 2
$ws.Range("E16").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -50
$ws.Range("E16").Copy($ws.Range("H22"))
$ws.Range("H22").Value = -50

# --- Plain numeric value updates ---
$ws.Range("M15").Value = 0
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 40
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 40
$ws.Range("I16").Value = 185
$ws.Range("J16").Value = 121
$ws.Range("K16").Value = 52.892561983471
$ws.Range("L16").Value = 12.121212121212
$ws.Range("M16").Value = -22.268907563025
$ws.Range("N16").Value = -79.143179255918
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 133.333333333333
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = 31.818181818181
$ws.Range("I17").Value = 274
$ws.Range("J17").Value = 229
$ws.Range("K17").Value = 19.650655021834
$ws.Range("L17").Value = 3.787878787878
$ws.Range("M17").Value = 39.795918367346
$ws.Range("N17").Value = 1.481481481481
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -42.857142857142
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 28
$ws.Range("H18").Value = -53.571428571428
$ws.Range("I18").Value = 240
$ws.Range("J18").Value = 205
$ws.Range("K18").Value = 17.073170731707
$ws.Range("L18").Value = -16.083916083916
$ws.Range("M18").Value = -40.446650124069
$ws.Range("N18").Value = -86.599664991624
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -8.333333333333
$ws.Range("F19").Value = 47
$ws.Range("G19").Value = 45
$ws.Range("H19").Value = 4.444444444444
$ws.Range("I19").Value = 577
$ws.Range("J19").Value = 479
$ws.Range("K19").Value = 20.459290187891
$ws.Range("L19").Value = 8.867924528301
$ws.Range("M19").Value = 50.260416666666
$ws.Range("N19").Value = -1.198630136986
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 29
$ws.Range("H20").Value = 31.818181818181
$ws.Range("I20").Value = 279
$ws.Range("J20").Value = 192
$ws.Range("K20").Value = 45.3125
$ws.Range("L20").Value = 74.375
$ws.Range("M20").Value = -14.678899082568
$ws.Range("N20").Value = -91.103316326530
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = 21.875
$ws.Range("F21").Value = 140
$ws.Range("H21").Value = 6.060606060606
$ws.Range("I21").Value = 1576
$ws.Range("J21").Value = 1240
$ws.Range("K21").Value = 27.096774193548
$ws.Range("L21").Value = 11.220889202540
$ws.Range("M21").Value = 0.318268618714
$ws.Range("N21").Value = -76.495152870991
$ws.Range("I22").Value = 14
$ws.Range("J22").Value = 8
$ws.Range("K22").Value = 75
$ws.Range("L22").Value = -26.315789473684
$ws.Range("M22").Value = -22.222222222222
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -15.151515151515
$ws.Range("F24").Value = 119
$ws.Range("G24").Value = 109
$ws.Range("H24").Value = 9.174311926605
$ws.Range("I24").Value = 1312
$ws.Range("J24").Value = 1198
$ws.Range("K24").Value = 9.515859766277
$ws.Range("L24").Value = 0.382555470543
$ws.Range("M24").Value = 40.021344717182
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 71.428571428571
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = 40
$ws.Range("I25").Value = 482
$ws.Range("J25").Value = 446
$ws.Range("K25").Value = 8.071748878923
$ws.Range("L25").Value = 17.560975609756
$ws.Range("M25").Value = -23.974763406940
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 64
$ws.Range("K27").Value = 36.170212765957
$ws.Range("L27").Value = 42.222222222222
$ws.Range("I28").Value = 10
$ws.Range("K28").Value = 233.333333333333
$ws.Range("L28").Value = 66.666666666666
$ws.Range("M28").Value = 150
$ws.Range("N28").Value = -61.538461538461
$ws.Range("I29").Value = 9
$ws.Range("K29").Value = 200
$ws.Range("L29").Value = 125
$ws.Range("M29").Value = 200
$ws.Range("N29").Value = -60.869565217391
